$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9272846058831874
$ws.Range("C2").Value = 0.2634518976413744
$ws.Range("E2").Value = 0.6412964237335359
$ws.Range("F2").Value = 2.281229430714916
$ws.Range("G2").Value = 0.5092661030865457
$ws.Range("H2").Value = 0.6058954889206092
$ws.Range("I2").Value = 0.4041742183378254
$ws.Range("J2").Value = 0.04263032953159041
$ws.Range("N2").Value = 0.8910042122132182
$ws.Range("B3").Value = 0.8198046374234877
$ws.Range("C3").Value = 0.2297148884699425
$ws.Range("E3").Value = 0.6235082885638974
$ws.Range("F3").Value = 2.236345635763882
$ws.Range("G3").Value = 0.4982493348151849
$ws.Range("H3").Value = 0.6063680181053002
$ws.Range("I3").Value = 0.4077551735070664
$ws.Range("J3").Value = 0.04221842670643738
$ws.Range("N3").Value = 0.8993479844436933
$ws.Range("B4").Value = 0.7538440884654847
$ws.Range("C4").Value = 0.2089680199773056
$ws.Range("E4").Value = 0.6128722674480471
$ws.Range("F4").Value = 2.210316002449588
$ws.Range("G4").Value = 0.4920306353164818
$ws.Range("H4").Value = 0.6071051769158942
$ws.Range("I4").Value = 0.4103307938338645
$ws.Range("J4").Value = 0.04200650065330791
$ws.Range("N4").Value = 0.9049141731238635
$ws.Range("B5").Value = 0.7269734137487376
$ws.Range("C5").Value = 0.2005053538736377
$ws.Range("E5").Value = 0.6086098355900162
$ws.Range("F5").Value = 2.200091814189193
$ws.Range("G5").Value = 0.4896326738082593
$ws.Range("H5").Value = 0.6075176360250367
$ws.Range("I5").Value = 0.411474732027127
$ws.Range("J5").Value = 0.04193038097091417
$ws.Range("N5").Value = 0.9072938784523217
$ws.Range("B6").Value = 0.722512106970612
$ws.Range("C6").Value = 0.1990996411257697
$ws.Range("E6").Value = 0.6079063984795567
$ws.Range("F6").Value = 2.198417192211423
$ws.Range("G6").Value = 0.4892426917128603
$ws.Range("H6").Value = 0.6075928826437149
$ws.Range("I6").Value = 0.411670368380463
$ws.Range("J6").Value = 0.04191835807998245
$ws.Range("N6").Value = 0.9076957581016885
$ws.Range("B7").Value = 0.75348166478949
$ws.Range("C7").Value = 0.2088539225558463
$ws.Range("E7").Value = 0.6128144920489262
$ws.Range("F7").Value = 2.210176566217541
$ws.Range("G7").Value = 0.4919977453850208
$ws.Range("H7").Value = 0.6071102862168658
$ws.Range("I7").Value = 0.4103458399192448
$ws.Range("J7").Value = 0.04200543269177359
$ws.Range("N7").Value = 0.9049458154238295
$ws.Range("B8").Value = 0.8902190672761208
$ws.Range("C8").Value = 0.2518260207222056
$ws.Range("E8").Value = 0.6351036692010723
$ws.Range("F8").Value = 2.265435228476164
$ws.Range("G8").Value = 0.5053536681960935
$ws.Range("H8").Value = 0.6059654260014895
$ws.Range("I8").Value = 0.4053304513201006
$ws.Range("J8").Value = 0.04247976046560353
$ws.Range("N8").Value = 0.8937891953170549
$ws.Range("B9").Value = 1.158607299643222
$ws.Range("C9").Value = 0.3358451710723216
$ws.Range("E9").Value = 0.6810881866650362
$ws.Range("F9").Value = 2.386004436541981
$ws.Range("G9").Value = 0.535921121353411
$ws.Range("H9").Value = 0.6072838587067366
$ws.Range("I9").Value = 0.3985051194819036
$ws.Range("J9").Value = 0.04373810663318167
$ws.Range("N9").Value = 0.875426176943229
$ws.Range("B10").Value = 1.355947102412756
$ws.Range("C10").Value = 0.397437545182413
$ws.Range("E10").Value = 0.7162737084558159
$ws.Range("F10").Value = 2.482143533284869
$ws.Range("G10").Value = 0.5611153989909212
$ws.Range("H10").Value = 0.6104491749223939
$ws.Range("I10").Value = 0.3953524996929403
$ws.Range("J10").Value = 0.04486701852786723
$ws.Range("N10").Value = 0.8640776577463853
$ws.Range("B11").Value = 1.445758203071875
$ws.Range("C11").Value = 0.4254315965901014
$ws.Range("E11").Value = 0.732587914826567
$ws.Range("F11").Value = 2.527546458624442
$ws.Range("G11").Value = 0.5731860576235448
$ws.Range("H11").Value = 0.6123714998538787
$ws.Range("I11").Value = 0.3943284823118276
$ws.Range("J11").Value = 0.04542592343729979
$ws.Range("N11").Value = 0.8593803704753284
$ws.Range("B12").Value = 1.479772873619652
$ws.Range("C12").Value = 0.4360288444809157
$ws.Range("E12").Value = 0.7388101697025036
$ws.Range("F12").Value = 2.544981194471603
$ws.Range("G12").Value = 0.5778457256316898
$ws.Range("H12").Value = 0.6131692218982892
$ws.Range("I12").Value = 0.3940001676072313
$ws.Range("J12").Value = 0.04564416285742823
$ws.Range("N12").Value = 0.8576685564841497
$ws.Range("B13").Value = 1.472446986295665
$ws.Range("C13").Value = 0.4337466919985786
$ws.Range("E13").Value = 0.7374681184387981
$ws.Range("F13").Value = 2.541215538694985
$ws.Range("G13").Value = 0.5768382169797803
$ws.Range("H13").Value = 0.6129943081663498
$ws.Range("I13").Value = 0.3940682239054567
$ws.Range("J13").Value = 0.0455968666612776
$ws.Range("N13").Value = 0.8580342479042145
$ws.Range("B14").Value = 1.448556510718504
$ws.Range("C14").Value = 0.4263035082313991
$ws.Range("E14").Value = 0.7330989326114832
$ws.Range("F14").Value = 2.528975972526013
$ws.Range("G14").Value = 0.5735676258590416
$ws.Range("H14").Value = 0.6124357280857993
$ws.Range("I14").Value = 0.3943002772527961
$ws.Range("J14").Value = 0.04544374555347019
$ws.Range("N14").Value = 0.8592381964829059
$ws.Range("B15").Value = 1.433923557583796
$ws.Range("C15").Value = 0.421743890071582
$ws.Range("E15").Value = 0.7304284681367932
$ws.Range("F15").Value = 2.521510402438196
$ws.Range("G15").Value = 0.5715758878897503
$ws.Range("H15").Value = 0.6121026812042629
$ws.Range("I15").Value = 0.3944501742482558
$ws.Range("J15").Value = 0.04535081536453589
$ws.Range("N15").Value = 0.8599843698039535
$ws.Range("B16").Value = 1.350078491257818
$ws.Range("C16").Value = 0.3956075766478762
$ws.Range("E16").Value = 0.7152137464171204
$ws.Range("F16").Value = 2.47921006722791
$ws.Range("G16").Value = 0.5603389151255413
$ws.Range("H16").Value = 0.6103332886791719
$ws.Range("I16").Value = 0.3954277156335237
$ws.Range("J16").Value = 0.04483141162020132
$ws.Range("N16").Value = 0.8643940018529221
$ws.Range("B17").Value = 1.298652111627575
$ws.Range("C17").Value = 0.3795675209228193
$ws.Range("E17").Value = 0.7059590046371085
$ws.Range("F17").Value = 2.453688772807197
$ws.Range("G17").Value = 0.5536023274106157
$ws.Range("H17").Value = 0.6093716701516314
$ws.Range("I17").Value = 0.3961327913907482
$ws.Range("J17").Value = 0.04452444236100916
$ws.Range("N17").Value = 0.8672183531916389
$ws.Range("B18").Value = 1.269076842146319
$ws.Range("C18").Value = 0.3703393837253657
$ws.Range("E18").Value = 0.7006649183010438
$ws.Range("F18").Value = 2.439166570730492
$ws.Range("G18").Value = 0.5497849376586856
$ws.Range("H18").Value = 0.6088639542845584
$ws.Range("I18").Value = 0.3965769062084377
$ws.Range("J18").Value = 0.04435214667101661
$ws.Range("N18").Value = 0.8688866297626987
$ws.Range("B19").Value = 1.259063840235171
$ws.Range("C19").Value = 0.3672144948831146
$ws.Range("E19").Value = 0.6988774085377543
$ws.Range("F19").Value = 2.434276512188859
$ws.Range("G19").Value = 0.548502245842144
$ws.Range("H19").Value = 0.6086998323368675
$ws.Range("I19").Value = 0.3967338862857446
$ws.Range("J19").Value = 0.04429454062800176
$ws.Range("N19").Value = 0.8694589981607805
$ws.Range("B20").Value = 1.304126146376746
$ws.Range("C20").Value = 0.3812752521506582
$ws.Range("E20").Value = 0.706941185910182
$ws.Range("F20").Value = 2.456389300189898
$ws.Range("G20").Value = 0.5543135085636806
$ws.Range("H20").Value = 0.6094693362884982
$ws.Range("I20").Value = 0.3960537391252856
$ws.Range("J20").Value = 0.04455667784208828
$ws.Range("N20").Value = 0.8669131643068582
$ws.Range("B21").Value = 1.455573587534104
$ws.Range("C21").Value = 0.4284898460776958
$ws.Range("E21").Value = 0.7343810619825604
$ws.Range("F21").Value = 2.53256445855888
$ws.Range("G21").Value = 0.5745258596931109
$ws.Range("H21").Value = 0.6125978995221004
$ws.Range("I21").Value = 0.3942305000131583
$ws.Range("J21").Value = 0.04548854139208913
$ws.Range("N21").Value = 0.8588827502499399
$ws.Range("B22").Value = 1.554583014680475
$ws.Range("C22").Value = 0.4593269858704048
$ws.Range("E22").Value = 0.7525735697339542
$ws.Range("F22").Value = 2.583758540159863
$ws.Range("G22").Value = 0.5882537442579121
$ws.Range("H22").Value = 0.6150495162920038
$ws.Range("I22").Value = 0.3933857030779606
$ws.Range("J22").Value = 0.04613603822117085
$ws.Range("N22").Value = 0.8540246558138946
$ws.Range("B23").Value = 1.501737300845718
$ws.Range("C23").Value = 0.4428704662896052
$ws.Range("E23").Value = 0.74284015332249
$ws.Range("F23").Value = 2.556305796489283
$ws.Range("G23").Value = 0.5808791558105781
$ws.Range("H23").Value = 0.6137036728115106
$ws.Range("I23").Value = 0.3938046951200604
$ws.Range("J23").Value = 0.04578691305971461
$ws.Range("N23").Value = 0.8565817863363208
$ws.Range("B24").Value = 1.301651364997383
$ws.Range("C24").Value = 0.3805032071078926
$ws.Range("E24").Value = 0.7064970589479884
$ws.Range("F24").Value = 2.455167923601891
$ws.Range("G24").Value = 0.5539918105656767
$ws.Range("H24").Value = 0.6094250408674782
$ws.Range("I24").Value = 0.3960893579578268
$ws.Range("J24").Value = 0.04454209115192498
$ws.Range("N24").Value = 0.8670510015203519
$ws.Range("B25").Value = 1.085974339969994
$ws.Range("C25").Value = 0.3131408562290972
$ws.Range("E25").Value = 0.6684031713612484
$ws.Range("F25").Value = 2.352068583461232
$ws.Range("G25").Value = 0.5271762343753181
$ws.Range("H25").Value = 0.6065431949861448
$ws.Range("I25").Value = 0.4000265391705327
$ws.Range("J25").Value = 0.04336210953906061
$ws.Range("N25").Value = 0.8800175854510357

Write-Host "Applied 216 cell updates"
